$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.790.07'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.649.04'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.66%  '
$ws.Range('D5').Value = "'216.87"
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').Value = "'19.28"
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').Value = "'0.0845"
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '1.872.92'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'4.22"
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.627.23'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').Value = "'0.533"
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = "'65.70"
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '26.791.31'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = "'217.20"
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').Value = "'4.38"
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').Value = "'2.44"
$ws.Range('E22').Value = '  +16.30%  '
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = "'145.74"
$ws.Range('E25').Value = '  -1.43%  '
$ws.Range('E26').Value = '  +0.47%  '
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').Value = "'7.22"
$ws.Range('E28').Value = '  +4.01%  '
$ws.Range('D29').Value = "'15.82"
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.278.27'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'1.55"
$ws.Range('E35').Value = '  +0.85%  '
$ws.Range('D36').Value = "'2.44"
$ws.Range('E36').Value = '  +2.15%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').Value = '  +5.47%  '
$ws.Range('E39').Value = '  +3.12%  '
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('D41').Value = "'0.819"
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('E42').Value = '  -1.53%  '
$ws.Range('E43').Value = '  +1.23%  '
$ws.Range('D44').Value = '1.798.17'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = "'92.21"
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('D46').Value = "'59.60"
$ws.Range('E46').Value = '  +6.51%  '
$ws.Range('D47').Value = "'1.63"
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = "'7.79"
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('D51').Value = "'0.0981"
$ws.Range('E51').Value = '  +1.41%  '
